$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 0.06427130384423536
$ws.Range("H2").Value2 = -0.04804879519399039
$ws.Range("I2").Value2 = -23.60965181406512
$ws.Range("G3").Value2 = 0.06437240657180693
$ws.Range("H3").Value2 = 14.52730255528484
$ws.Range("G4").Value2 = -0.01406049033992971
$ws.Range("H4").Value2 = 30.77608474799679
$ws.Range("G5").Value2 = -0.02125627512239815
$ws.Range("H5").Value2 = -85.94750855186763
$ws.Range("G6").Value2 = -0.02470200997972323
$ws.Range("H6").Value2 = -120.6805978042984
$ws.Range("G7").Value2 = -0.0235241678783848
$ws.Range("H7").Value2 = -315.3471812927345
$ws.Range("G8").Value2 = -0.001765022885003737
$ws.Range("H8").Value2 = 68.97570194638635
$ws.Range("G9").Value2 = -0.001192280637968696
$ws.Range("H9").Value2 = 78.29512293277574
$ws.Range("G10").Value2 = -0.0626825797351365
$ws.Range("H10").Value2 = 0.5381870428465005
$ws.Range("G11").Value2 = -0.06976873356215615
$ws.Range("H11").Value2 = -8.88086599560223
$ws.Range("G12").Value2 = -0.4076152984582946
$ws.Range("H12").Value2 = -3.297462105419792
$ws.Range("G13").Value2 = -0.4005406476591079
$ws.Range("H13").Value2 = -2.147675379258007
$ws.Range("G14").Value2 = -0.01890752586753613
$ws.Range("H14").Value2 = -130.7410035261609
$ws.Range("G15").Value2 = -0.01562181837469445
$ws.Range("H15").Value2 = 65.52896081650881
$ws.Range("G16").Value2 = 0.1305660323562913
$ws.Range("H16").Value2 = -4.521852165802985
$ws.Range("G17").Value2 = 0.1425948927336241
$ws.Range("H17").Value2 = 2.24008024258985
$ws.Range("G18").Value2 = 0.1215320387719089
$ws.Range("H18").Value2 = 3.283666013587971
$ws.Range("G19").Value2 = 0.1167391776739392
$ws.Range("H19").Value2 = -9.271609598883257
$ws.Range("G20").Value2 = 0.08459635290887559
$ws.Range("H20").Value2 = -4.663657988381644
$ws.Range("G21").Value2 = 0.09238025256256767
$ws.Range("H21").Value2 = 6.073217394008084
$ws.Range("G22").Value2 = -0.1033337862835137
$ws.Range("H22").Value2 = -10.53731590326026
$ws.Range("G23").Value2 = -0.09929440307253601
$ws.Range("H23").Value2 = 2.119980942921455
$ws.Range("G24").Value2 = 0.1725305686686509
$ws.Range("H24").Value2 = 7.101267491811485
$ws.Range("G25").Value2 = 0.1644516245163774
$ws.Range("H25").Value2 = -3.601559884240534
$ws.Range("G26").Value2 = 0.08889646714925685
$ws.Range("H26").Value2 = -1.938505195803601
$ws.Range("G27").Value2 = 0.08466480802646421
$ws.Range("H27").Value2 = -1.508738839377272
$ws.Range("G28").Value2 = -0.1448863674616052
$ws.Range("H28").Value2 = -5.264391258886106
$ws.Range("G29").Value2 = -0.1280315681962876
$ws.Range("H29").Value2 = 8.417894860547534
$ws.Range("G30").Value2 = 0.04647237872779753
$ws.Range("H30").Value2 = -10.65628491920817
$ws.Range("G31").Value2 = 0.04215354341563844
$ws.Range("H31").Value2 = -3.801648317328129
$ws.Range("G32").Value2 = 0.1159157970289918
$ws.Range("H32").Value2 = 6.634410782776333
$ws.Range("G33").Value2 = 0.1207633875778223
$ws.Range("H33").Value2 = -2.684878273361186
$ws.Range("G34").Value2 = -0.01067907096475751
$ws.Range("H34").Value2 = 31.61628093517611
$ws.Range("G35").Value2 = -0.01370548079639382
$ws.Range("H35").Value2 = 18.11409074756326
$ws.Range("G36").Value2 = 0.03295728617703193
$ws.Range("H36").Value2 = -10.36155443677626
$ws.Range("G37").Value2 = 0.03043665691321512
$ws.Range("H37").Value2 = -14.7089706255889
$ws.Range("G38").Value2 = 0.09776095486946634
$ws.Range("H38").Value2 = -2.533617170116092
$ws.Range("G39").Value2 = 0.1038446446876256
$ws.Range("H39").Value2 = 6.616661938911427
$ws.Range("G40").Value2 = 0.03031990966901548
$ws.Range("H40").Value2 = -9.99687749786265
$ws.Range("G41").Value2 = 0.03456963301660645
$ws.Range("H41").Value2 = 7.296713261453721
$ws.Range("G42").Value2 = 0.1190174342368622
$ws.Range("H42").Value2 = -1.561547931283802
$ws.Range("G43").Value2 = 0.126764968494164
$ws.Range("H43").Value2 = -0.7985875843607912
$ws.Range("G44").Value2 = 0.03438256321436341
$ws.Range("H44").Value2 = -13.3212092233607
$ws.Range("G45").Value2 = 0.03870755936427836
$ws.Range("H45").Value2 = 24.21256973029249
$ws.Range("G46").Value2 = 0.06057773736683784
$ws.Range("H46").Value2 = 6.994411294244696
$ws.Range("G47").Value2 = 0.05884383917028591
$ws.Range("H47").Value2 = 0.2959808337499061
$ws.Range("G48").Value2 = 0.04146497924594267
$ws.Range("H48").Value2 = -15.81966506863696
$ws.Range("G49").Value2 = 0.03941549345120769
$ws.Range("H49").Value2 = -13.50133416001853
$ws.Range("G50").Value2 = 0.02942323670676247
$ws.Range("H50").Value2 = 11.08774188579917
$ws.Range("G51").Value2 = 0.02942323670676247
$ws.Range("H51").Value2 = 5.023468213285062
$ws.Range("G52").Value2 = -0.08420345248284296
$ws.Range("H52").Value2 = 3.12940023394788
$ws.Range("G53").Value2 = -0.07873385448778604
$ws.Range("H53").Value2 = 1.847940198659941
$ws.Range("G54").Value2 = 0.05276975029654993
$ws.Range("H54").Value2 = 5.4831685151142
$ws.Range("G55").Value2 = 0.04874859276054085
$ws.Range("H55").Value2 = -13.389813297607
$ws.Range("G56").Value2 = 0.04511370100157242
$ws.Range("H56").Value2 = -8.73564322572272
$ws.Range("G57").Value2 = 0.0392146256705742
$ws.Range("H57").Value2 = 3.239912936494189
$ws.Range("G58").Value2 = 0.0536721592442644
$ws.Range("H58").Value2 = -6.841895660731708
$ws.Range("G59").Value2 = 0.06062573782934053
$ws.Range("H59").Value2 = 6.336114484101742
$ws.Range("G60").Value2 = 0.02732655038723472
$ws.Range("H60").Value2 = -0.4859764514020445
$ws.Range("G61").Value2 = 0.03175514860336491
$ws.Range("H61").Value2 = 18.94282858819619
$ws.Range("G62").Value2 = 0.06633489006864447
$ws.Range("H62").Value2 = 6.225510270950167
$ws.Range("G63").Value2 = 0.06746974183425647
$ws.Range("H63").Value2 = 5.615094854880258
$ws.Range("G64").Value2 = 0.03207501596671017
$ws.Range("H64").Value2 = 15.62873858474947
$ws.Range("G65").Value2 = 0.03191452030104739
$ws.Range("H65").Value2 = -9.91533754919374
$ws.Range("G66").Value2 = 0.08480380953122589
$ws.Range("H66").Value2 = 9.166241982503157
$ws.Range("G67").Value2 = 0.07889599869896183
$ws.Range("H67").Value2 = 0.03190729298184846
$ws.Range("G68").Value2 = -0.01822101453720723
$ws.Range("H68").Value2 = 16.20862687354899
$ws.Range("G69").Value2 = -0.0162609380114367
$ws.Range("H69").Value2 = 15.04794527876746
$ws.Range("G70").Value2 = 0.06468597564025703
$ws.Range("H70").Value2 = -10.14072333415797
$ws.Range("G71").Value2 = 0.07458972733623058
$ws.Range("H71").Value2 = -6.08316214604926
$ws.Range("G72").Value2 = -0.1430387745441365
$ws.Range("H72").Value2 = 6.900940843926855
$ws.Range("G73").Value2 = -0.1502113743438646
$ws.Range("H73").Value2 = 1.87108940726394
$ws.Range("G74").Value2 = 0.1543301687381017
$ws.Range("H74").Value2 = 2.600050410734917
$ws.Range("G75").Value2 = 0.1561595334246783
$ws.Range("H75").Value2 = 3.796816529111589
$ws.Range("G76").Value2 = -0.01534922450543432
$ws.Range("H76").Value2 = -1380.925450794978
$ws.Range("G77").Value2 = -0.008570633565215304
$ws.Range("H77").Value2 = -288.1824802599772
$ws.Range("G78").Value2 = 0.1000163909531586
$ws.Range("H78").Value2 = 11.18089835009147
$ws.Range("G79").Value2 = 0.09335006138216749
$ws.Range("H79").Value2 = -3.664864115276225
$ws.Range("G80").Value2 = -0.220584731623212
$ws.Range("H80").Value2 = -1.918079360261572
$ws.Range("G81").Value2 = -0.2101152744319672
$ws.Range("H81").Value2 = 1.407046241060185
$ws.Range("G82").Value2 = 0.1708251674848299
$ws.Range("H82").Value2 = 1.913580659210046
$ws.Range("G83").Value2 = 0.1739798400300795
$ws.Range("H83").Value2 = -1.164065548307729
$ws.Range("G84").Value2 = 0.1088098352863314
$ws.Range("H84").Value2 = 2.525342431506229
$ws.Range("G85").Value2 = 0.1165955093496981
$ws.Range("H85").Value2 = 11.51189727695051
